# Reset COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE counters (columns D, E)
# and the FORMATCONSISTENCY / FORMATCONSISTENCY SCORE counters (columns T, U)
# for row 3 back to zero, reflecting the cleaned-up data processing pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8: zero out columns D (COMPLETENESSOPTIONAL) and E (COMPLETENESSOPTIONAL SCORE)
$ws.Range("D2:E8").Value = 0

# Row 3: zero out columns T (FORMATCONSISTENCY) and U (FORMATCONSISTENCY SCORE)
$ws.Range("T3:U3").Value = 0
